$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row Right value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row Right value (B12): 75 -> 125
$ws.Range("B12").Value = 125

# Update "Total" row Max text (E12): 74/84 -> 125/140
$ws.Range("E12").Value = "125/140"
